$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the contents of row 8 (A8:K8) while keeping their existing formatting.
$ws.Range("A8:K8").ClearContents()

# L34, L36, L37 and L38 were using a style that only differed from the
# regular column-L number style by an (invisible, no-color) fill flag.
# Dropping the explicit "no fill" pattern re-resolves them onto the plain
# number-format style used by the rest of the column.
$ws.Range("L34").Interior.Pattern = -4142  # xlNone
$ws.Range("L36").Interior.Pattern = -4142  # xlNone
$ws.Range("L37").Interior.Pattern = -4142  # xlNone
$ws.Range("L38").Interior.Pattern = -4142  # xlNone

# Update the view: select A8:L8 (the cleared row) as the active selection,
# matching the saved sheet view state.
$ws.Range("A8:L8").Select()
